# Update crypto price/volume data per latest scrape (GitHub Actions job)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.296.50"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.91%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.749.64"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.73%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.00"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "  +0.17%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "594.75"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.30%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.67"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.47%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.748.65"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -0.85%  "
$ws.Range("E8").Value = "  -0.05%  "
$ws.Range("E9").Value = "  -0.78%  "
$ws.Range("E10").Value = "  -3.42%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.03%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.448"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "  -1.15%  "
$ws.Range("E13").Value = "  -7.12%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.52%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.377.63"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.70%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.731.12"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.86%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.290.74"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.95%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.84"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -4.03%  "
$ws.Range("E19").Value = "  -2.53%  "
$ws.Range("E20").Value = "  -0.31%  "
$ws.Range("E21").Value = "  +1.37%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.10"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.41%  "
$ws.Range("E23").Value = "  -2.90%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.40"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.71%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.0000143"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -2.99%  "
$ws.Range("E26").Value = "  -0.81%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "11.98"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  -1.64%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.45%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.896.35"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.69%  "
$ws.Range("E31").Value = "  -4.64%  "
$ws.Range("E32").Value = "  -4.61%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.85"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.25%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.29%  "
$ws.Range("E35").Value = "  +0.81%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.704.99"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -0.82%  "
$ws.Range("E38").Value = "  -2.88%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.38"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.19%  "
$ws.Range("E40").Value = "  +0.30%  "
$ws.Range("E41").Value = "  -0.42%  "
$ws.Range("E42").Value = "  -0.03%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.09%  "
$ws.Range("E44").Value = "  +0.00%  "
$ws.Range("E45").Value = "  -2.39%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.59"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.89%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "1.92"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.26%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.82"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.82"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.09%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.04"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +3.94%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "390.37"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.25%  "
